# Generate Report for Handback
# Marks the zh-cn / de-de handback rows as complete: fills in the
# "Latest Target File" (hyperlink) + "Latest Handback File" + "Latest
# Handback DateTime" columns, updates the Status text, and widens a few
# columns to fit the new (longer) values.

$wb = $excel.ActiveWorkbook

# ---- Overview sheet: widen the zh-cn / de-de status columns (E, F) ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Columns.Item(5).ColumnWidth = 29.17
$wsOverview.Columns.Item(6).ColumnWidth = 29.17
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"

# ---- zh-cn sheet ----
$wsZhCn = $wb.Worksheets.Item("zh-cn")

# Status column (C) widened + text updated
$wsZhCn.Columns.Item(3).ColumnWidth = 29.17
$wsZhCn.Range("C2").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("C3").Value = "Handed back: in sync with en-US"

# Latest Target File / Latest Handback File columns (I, J) widened
$wsZhCn.Columns.Item(9).ColumnWidth = 39.17
$wsZhCn.Columns.Item(10).ColumnWidth = 39.17

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/94bf12ac779f9195801d3063c8319ee8c161abf9/e2e/050d8661-4033-43ea-927c-368262d5cc9a.md", "", "", "050d8661-4033-43ea-927c-368262d5cc9a.md")
$wsZhCn.Range("J2").Value = "050d8661-4033-43ea-927c-368262d5cc9a.1eae97f96c692ca7a407afb985f23505896d36e5.zh-cn.xlf"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/94bf12ac779f9195801d3063c8319ee8c161abf9/e2e/333f9a0b-1b5a-4f41-97af-98972a2ccf4b.md", "", "", "333f9a0b-1b5a-4f41-97af-98972a2ccf4b.md")
$wsZhCn.Range("J3").Value = "333f9a0b-1b5a-4f41-97af-98972a2ccf4b.37127dacff6824b53ec2418edd81cd882d658e2b.zh-cn.xlf"

# Latest Handback DateTime (K) - zh-cn finished first
$wsZhCn.Range("K2").Value = "2016-08-21 12:37:16"
$wsZhCn.Range("K3").Value = "2016-08-21 12:37:16"

# ---- de-de sheet ----
$wsDeDe = $wb.Worksheets.Item("de-de")

# Status column (C) widened + text updated
$wsDeDe.Columns.Item(3).ColumnWidth = 29.17
$wsDeDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("C3").Value = "Handed back: in sync with en-US"

# Latest Target File / Latest Handback File columns (I, J) widened
$wsDeDe.Columns.Item(9).ColumnWidth = 39.17
$wsDeDe.Columns.Item(10).ColumnWidth = 39.17

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/94bf12ac779f9195801d3063c8319ee8c161abf9/e2e/050d8661-4033-43ea-927c-368262d5cc9a.md", "", "", "050d8661-4033-43ea-927c-368262d5cc9a.md")
$wsDeDe.Range("J2").Value = "050d8661-4033-43ea-927c-368262d5cc9a.1eae97f96c692ca7a407afb985f23505896d36e5.de-de.xlf"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/94bf12ac779f9195801d3063c8319ee8c161abf9/e2e/333f9a0b-1b5a-4f41-97af-98972a2ccf4b.md", "", "", "333f9a0b-1b5a-4f41-97af-98972a2ccf4b.md")
$wsDeDe.Range("J3").Value = "333f9a0b-1b5a-4f41-97af-98972a2ccf4b.37127dacff6824b53ec2418edd81cd882d658e2b.de-de.xlf"

# Latest Handback DateTime (K) - de-de finished a few seconds later
$wsDeDe.Range("K2").Value = "2016-08-21 12:37:16"
$wsDeDe.Range("K3").Value = "2016-08-21 12:37:22"

Write-Output "Handback report generated"
